$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 10

$b64 = "VG8gZmluZCB0aGUgdGhyZWUtZGlnaXQgbnVtYmVyLCBsZXQncyBkZW5vdGUgdGhlIGRpZ2l0cyBhcyBmb2xsb3dzOgotIExldCB0aGUgZmlyc3QgZGlnaXQgYmUgXCggYSBcKS4KLSBMZXQgdGhlIHNlY29uZCBkaWdpdCBiZSBcKCBiIFwpLgotIExldCB0aGUgdGhpcmQgZGlnaXQgYmUgXCggYyBcKS4KCkZyb20gdGhlIHByb2JsZW0sIHdlIGhhdmUgdGhlIGZvbGxvd2luZyByZWxhdGlvbnNoaXBzOgoxLiBUaGUgc2Vjb25kIGRpZ2l0IFwoIGIgXCkgaXMgZm91ciB0aW1lcyB0aGUgdGhpcmQgZGlnaXQgXCggYyBcKToKICAgXFsKICAgYiA9IDRjCiAgIFxdCgoyLiBUaGUgZmlyc3QgZGlnaXQgXCggYSBcKSBpcyB0aHJlZSBsZXNzIHRoYW4gdGhlIHNlY29uZCBkaWdpdCBcKCBiIFwpOgogICBcWwogICBhID0gYiAtIDMKICAgXF0KClNpbmNlIFwoIGEgXCksIFwoIGIgXCksIGFuZCBcKCBjIFwpIGFyZSBkaWdpdHMsIHRoZXkgbXVzdCBiZSBpbnRlZ2VycyBiZXR3ZWVuIDAgYW5kIDkgaW5jbHVzaXZlLiBMZXQncyBhbmFseXplIHRoZSBjb25zdHJhaW50cyBzdGVwIGJ5IHN0ZXAuCgpGaXJzdCwgc2luY2UgXCggYiA9IDRjIFwpLCBcKCBiIFwpIG11c3QgYmUgYSBtdWx0aXBsZSBvZiA0LiBUaGUgcG9zc2libGUgdmFsdWVzIGZvciBcKCBjIFwpIChzaW5jZSBcKCBjIFwpIGlzIGEgZGlnaXQpIGFyZToKXFsKYyA9IDAsIDEsIDIKXF0KCkxldCdzIGV2YWx1YXRlIGVhY2ggY2FzZToKCjEuICoqSWYgXCggYyA9IDAgXCkqKjoKICAgXFsKICAgYiA9IDQgXHRpbWVzIDAgPSAwCiAgIFxdCiAgIFxbCiAgIGEgPSAwIC0gMyA9IC0zCiAgIFxdCiAgIFRoaXMgaXMgbm90IHBvc3NpYmxlIHNpbmNlIFwoIGEgXCkgbXVzdCBiZSBhIGRpZ2l0IGJldHdlZW4gMCBhbmQgOS4KCjIuICoqSWYgXCggYyA9IDEgXCkqKjoKICAgXFsKICAgYiA9IDQgXHRpbWVzIDEgPSA0CiAgIFxdCiAgIFxbCiAgIGEgPSA0IC0gMyA9IDEKICAgXF0KICAgSGVyZSwgXCggYSA9IDEgXCksIFwoIGIgPSA0IFwpLCBhbmQgXCggYyA9IDEgXCkuIFRoaXMgaXMgYSB2YWxpZCBzZXQgb2YgZGlnaXRzLgoKMy4gKipJZiBcKCBjID0gMiBcKSoqOgogICBcWwogICBiID0gNCBcdGltZXMgMiA9IDgKICAgXF0KICAgXFsKICAgYSA9IDggLSAzID0gNQogICBcXQogICBIZXJlLCBcKCBhID0gNSBcKSwgXCggYiA9IDggXCksIGFuZCBcKCBjID0gMiBcKS4gVGhpcyBpcyBhbHNvIGEgdmFsaWQgc2V0IG9mIGRpZ2l0cy4KClRodXMsIHdlIGhhdmUgdHdvIHBvc3NpYmxlIHRocmVlLWRpZ2l0IG51bWJlcnM6Ci0gXCggMTQxIFwpIChmcm9tIFwoIGEgPSAxIFwpLCBcKCBiID0gNCBcKSwgXCggYyA9IDEgXCkpCi0gXCggNTgyIFwpIChmcm9tIFwoIGEgPSA1IFwpLCBcKCBiID0gOCBcKSwgXCggYyA9IDIgXCkpCgpIb3dldmVyLCB0aGUgcHJvYmxlbSBzcGVjaWZpZXMgdGhhdCB0aGUgc2Vjb25kIGRpZ2l0IGlzIGZvdXIgdGltZXMgdGhlIHRoaXJkIGRpZ2l0LCBhbmQgdGhlIGZpcnN0IGRpZ2l0IGlzIHRocmVlIGxlc3MgdGhhbiB0aGUgc2Vjb25kIGRpZ2l0LiBCb3RoIGNvbmRpdGlvbnMgYXJlIHNhdGlzZmllZCBpbiBib3RoIGNhc2VzLiBUaGVyZWZvcmUsIHRoZSB0aHJlZS1kaWdpdCBudW1iZXIgY291bGQgYmUgZWl0aGVyIFwoIDE0MSBcKSBvciBcKCA1ODIgXCku"
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("C2").Value = $text

$b64 = "VGhlIG1vZGVsIGNvcnJlY3RseSBpZGVudGlmaWVkIGJvdGggcG9zc2libGUgdGhyZWUtZGlnaXQgbnVtYmVycywgMTQxIGFuZCA1ODIsIGFzIGV4cGVjdGVkLg=="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("E2").Value = $text

$b64 = "TGV0J3MgYnJlYWsgdGhpcyBkb3duIHN0ZXAgYnkgc3RlcDoKCjEuIFlvdSBjdXJyZW50bHkgaGF2ZSAzIGFwcGxlcy4KMi4gVGhlIGZhY3QgdGhhdCB5b3UgYXRlIG9uZSBhcHBsZSB5ZXN0ZXJkYXkgZG9lcyBub3QgYWZmZWN0IHRoZSBudW1iZXIgb2YgYXBwbGVzIHlvdSBoYXZlIHRvZGF5LgoKVGhlcmVmb3JlLCB5b3Ugc3RpbGwgaGF2ZSAzIGFwcGxlcyB0b2RheS4="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("C3").Value = $text

$b64 = "VGhlIG1vZGVsJ3MgYW5zd2VyIGNvcnJlY3RseSBleHBsYWlucyB3aHkgeW91IHN0aWxsIGhhdmUgMyBhcHBsZXMgdG9kYXku"
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("E3").Value = $text

$b64 = "VG8gZGV0ZXJtaW5lIGhvdyBsb25nIGl0IHdpbGwgdGFrZSB0byBkcnkgMjAgdG93ZWxzLCB3ZSBuZWVkIHRvIHVuZGVyc3RhbmQgdGhlIHJlbGF0aW9uc2hpcCBiZXR3ZWVuIHRoZSBudW1iZXIgb2YgdG93ZWxzIGFuZCB0aGUgZHJ5aW5nIHRpbWUuIAoKR2l2ZW46Ci0gSXQgdGFrZXMgMSBob3VyIHRvIGRyeSAxNSB0b3dlbHMuCgpBc3N1bWluZyB0aGUgZHJ5aW5nIHByb2Nlc3MgaXMgbGluZWFyIGFuZCB0aGUgZHJ5aW5nIGNhcGFjaXR5IGlzIGNvbnN0YW50LCB3ZSBjYW4gc2V0IHVwIGEgcHJvcG9ydGlvbiB0byBmaW5kIHRoZSB0aW1lIHJlcXVpcmVkIHRvIGRyeSAyMCB0b3dlbHMuCgpMZXQgXCggdCBcKSBiZSB0aGUgdGltZSBpbiBob3VycyB0byBkcnkgMjAgdG93ZWxzLgoKVGhlIHByb3BvcnRpb24gY2FuIGJlIHNldCB1cCBhcyBmb2xsb3dzOgpcWwpcZnJhY3sxIFx0ZXh0eyBob3VyfX17MTUgXHRleHR7IHRvd2Vsc319ID0gXGZyYWN7dCBcdGV4dHsgaG91cnN9fXsyMCBcdGV4dHsgdG93ZWxzfX0KXF0KClRvIHNvbHZlIGZvciBcKCB0IFwpLCB3ZSBjcm9zcy1tdWx0aXBseToKXFsKMSBcdGltZXMgMjAgPSAxNSBcdGltZXMgdApcXQoKVGhpcyBzaW1wbGlmaWVzIHRvOgpcWwoyMCA9IDE1dApcXQoKTmV4dCwgd2Ugc29sdmUgZm9yIFwoIHQgXCkgYnkgZGl2aWRpbmcgYm90aCBzaWRlcyBieSAxNToKXFsKdCA9IFxmcmFjezIwfXsxNX0gPSBcZnJhY3s0fXszfSBcdGV4dHsgaG91cnN9ClxdCgpUaHVzLCBpdCB3aWxsIHRha2UgXChcZnJhY3s0fXszfVwpIGhvdXJzLCBvciAxIGhvdXIgYW5kIDIwIG1pbnV0ZXMsIHRvIGRyeSAyMCB0b3dlbHMu"
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("C4").Value = $text

$b64 = "VGhlIG1vZGVsJ3MgYW5zd2VyIG9ubHkgYWRkcmVzc2VzIHRoZSBkcnlpbmcgdGltZSBmb3Igb25lIHNjZW5hcmlvIChkcnlpbmcgb25lIGJ5IG9uZSkgYW5kIGRvZXMgbm90IGNvbnNpZGVyIHRoZSBwYXJhbGxlbCBkcnlpbmcgY2FzZS4="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("E4").Value = $text

$b64 = "VG8gZGV0ZXJtaW5lIGhvdyBtYW55IHNpc3RlcnMgZWFjaCBvZiBKZXNzaWNhJ3MgYnJvdGhlcnMgaGF2ZSwgbGV0J3MgYnJlYWsgZG93biB0aGUgaW5mb3JtYXRpb24gZ2l2ZW46CgoxLiBKZXNzaWNhIGhhcyAyIGJyb3RoZXJzLgoyLiBKZXNzaWNhIGhhcyAxIHNpc3RlciAoaGVyc2VsZikuCgpOb3csIGxldCdzIGNvbnNpZGVyIHRoZSBwZXJzcGVjdGl2ZSBvZiBvbmUgb2YgSmVzc2ljYSdzIGJyb3RoZXJzOgotIEhlIGhhcyAxIHNpc3RlciAoSmVzc2ljYSkgYW5kIDEgYWRkaXRpb25hbCBzaXN0ZXIgKEplc3NpY2EncyBzaXN0ZXIpLgoKVGhlcmVmb3JlLCBlYWNoIG9mIEplc3NpY2EncyBicm90aGVycyBoYXMgMiBzaXN0ZXJzLgoKU28sIGVhY2ggb2YgSmVzc2ljYSdzIGJyb3RoZXJzIGhhcyAyIHNpc3RlcnMu"
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("C5").Value = $text

$b64 = "VGhlIG1vZGVsJ3MgYW5zd2VyIGNvcnJlY3RseSBpZGVudGlmaWVzIHRoYXQgZWFjaCBvZiBKZXNzaWNhJ3MgYnJvdGhlcnMgaGFzIDIgc2lzdGVycywgbWF0Y2hpbmcgdGhlIGV4cGVjdGVkIG91dHB1dC4="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("E5").Value = $text

$b64 = "VGhlIG1vZGVsJ3MgYW5zd2VyIGNvcnJlY3RseSBmb2xsb3dzIHRoZSBleHBlY3RlZCBvdXRwdXQncyBsb2dpYyBhbmQgc3RlcHMgdG8gZGV0ZXJtaW5lIHRoYXQgOS45IGlzIGdyZWF0ZXIgdGhhbiA5LjExLg=="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("E7").Value = $text

$b64 = "dGhlIHBhdHRlcm4gaXMgYWRkaW5nIHRoZW4gbWludXMgMS4gVGh1cywgMys0PTcsIDctMSA9IDYuIEZpbmFsIGFuc3dlcg=="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("B8").Value = $text

$b64 = "TGV0J3MgYW5hbHl6ZSB0aGUgZ2l2ZW4gZXF1YXRpb25zIHN0ZXAgYnkgc3RlcDoKCjEuIFRoZSBmaXJzdCBlcXVhdGlvbiBpcyBcKDEgKyAyID0gMlwpLiBOb3JtYWxseSwgXCgxICsgMlwpIGVxdWFscyBcKDNcKSwgYnV0IGhlcmUgaXQgaXMgZ2l2ZW4gYXMgXCgyXCkuIFRoaXMgc3VnZ2VzdHMgdGhhdCB0aGUgcmVzdWx0IG1pZ2h0IGJlIHRoZSBzZWNvbmQgbnVtYmVyIGluIHRoZSBlcXVhdGlvbi4KCjIuIFRoZSBzZWNvbmQgZXF1YXRpb24gaXMgXCgyICsgMyA9IDRcKS4gTm9ybWFsbHksIFwoMiArIDNcKSBlcXVhbHMgXCg1XCksIGJ1dCBoZXJlIGl0IGlzIGdpdmVuIGFzIFwoNFwpLiBBZ2FpbiwgdGhpcyBzdWdnZXN0cyB0aGF0IHRoZSByZXN1bHQgbWlnaHQgYmUgdGhlIHNlY29uZCBudW1iZXIgaW4gdGhlIGVxdWF0aW9uLgpGb2xsb3dpbmcgdGhpcyBwYXR0ZXJuLCB0aGUgcmVzdWx0IG9mIHRoZSBlcXVhdGlvbiBzZWVtcyB0byBiZSB0aGUgc2Vjb25kIG51bWJlciBpbiB0aGUgZXF1YXRpb24uCgozLiBUaGUgdGhpcmQgZXF1YXRpb24gaXMgXCgzICsgNCA9ID9cKS4gQWNjb3JkaW5nIHRvIHRoZSBwYXR0ZXJuIG9ic2VydmVkLCB0aGUgcmVzdWx0IHNob3VsZCBiZSB0aGUgc2Vjb25kIG51bWJlciBpbiB0aGUgZXF1YXRpb24sIHdoaWNoIGlzIFwoNFwpLgoKVGhlcmVmb3JlLCBcKDMgKyA0ID0gNFwpLg=="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("C8").Value = $text

$b64 = "VGhlIG1vZGVsJ3MgYW5zd2VyIGZvbGxvd3MgYSBkaWZmZXJlbnQgcGF0dGVybiBhbmQgYXJyaXZlcyBhdCBhbiBpbmNvcnJlY3QgZmluYWwgYW5zd2VyLg=="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("E8").Value = $text

$b64 = "MS4gU2hlIHJlYWNoZWQgaW50byB0aGUgYmFza2V0IGFuZCBwdWxsZWQgb3V0IGEgc2hpbnksIHJlZCBhcHBsZS4KMi4gVGhlIHRlYWNoZXIgcGxhY2VkIGEgZ29sZGVuIHN0YXIgc3RpY2tlciBvbiB0aGUgc3R1ZGVudCdzIGRyYXdpbmcgb2YgYW4gYXBwbGUuCjMuIEFmdGVyIGEgbG9uZyBoaWtlLCBoZSBzYXQgZG93biB1bmRlciBhIHRyZWUgYW5kIGVuam95ZWQgYSBjcmlzcCBhcHBsZS4KNC4gVGhlIHBpZSByZWNpcGUgY2FsbGVkIGZvciB0aHJlZSBjdXBzIG9mIHNsaWNlZCBhcHBsZS4KNS4gRHVyaW5nIHRoZSBwaWNuaWMsIGV2ZXJ5b25lIGVuam95ZWQgYSByZWZyZXNoaW5nIGJpdGUgb2YgYSBqdWljeSBhcHBsZS4KNi4gVGhlIGZhcm1lciBwcm91ZGx5IGRpc3BsYXllZCBoaXMgbGFyZ2VzdCBhbmQgbW9zdCBwZXJmZWN0IGFwcGxlLgo3LiBTaGUgZGVjaWRlZCB0byBtYWtlIGEgaGVhbHRoeSBzbmFjayBieSBjdXR0aW5nIHVwIGFuIGFwcGxlLgo4LiBUaGUgY2hpbGQgZWFnZXJseSBwb2ludGVkIHRvIHRoZSBmcnVpdCBzdGFuZCBhbmQgYXNrZWQgZm9yIGFuIGFwcGxlLgo5LiBIZSBwb2xpc2hlZCB0aGUgc3VyZmFjZSBvZiB0aGUgZGVzayBhbmQgcGxhY2VkIGEgc2luZ2xlLCBwZXJmZWN0IGFwcGxlLgoxMC4gVGhlIGFyb21hIG9mIGNpbm5hbW9uIGFuZCBiYWtlZCBhcHBsZSBmaWxsZWQgdGhlIGtpdGNoZW4u"
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("C9").Value = $text

$b64 = "QWxsIHNlbnRlbmNlcyBjb3JyZWN0bHkgZW5kIHdpdGggJ2FwcGxlJy4="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("E9").Value = $text

$b64 = "MS4gSSBoYXZlIG9uZSBhcHBsZSBpbiBteSBiYXNrZXQsIGFuZCBpdCdzIHRoZSBvbmx5IG9uZS4KMi4gU2hlIGJvdWdodCB0d28gdGlja2V0cyBmb3IgdGhlIGNvbmNlcnQsIHNvIHdlIGNhbiBib3RoIGdvLgozLiBUaGUgY2F0IGhhZCB0aHJlZSBraXR0ZW5zLCBhbGwgb2Ygd2hpY2ggd2VyZSBhZG9yYWJsZS4KNC4gSGUgZmluaXNoZWQgdGhlIHJhY2UgaW4gZm91cnRoIHBsYWNlLCB3aGljaCB3YXMgYSBncmVhdCBhY2hpZXZlbWVudC4KNS4gVGhlcmUgYXJlIGZpdmUgYm9va3Mgb24gdGhlIHNoZWxmLCBlYWNoIG9uZSBhIGRpZmZlcmVudCBnZW5yZS4KNi4gVGhlIHJlY2lwZSBjYWxscyBmb3Igc2l4IGVnZ3MsIHNvIG1ha2Ugc3VyZSB0byBidXkgZW5vdWdoLgo3LiBXZSBuZWVkIHNldmVuIHZvbHVudGVlcnMgdG8gaGVscCB3aXRoIHRoZSBldmVudCB0aGlzIHdlZWtlbmQuCjguIFRoZSBvY3RvcHVzIGhhcyBlaWdodCB0ZW50YWNsZXMsIGVhY2ggb25lIHdpdGggaXRzIG93biBmdW5jdGlvbi4KOS4gU2hlIGNlbGVicmF0ZWQgaGVyIG5pbnRoIGJpcnRoZGF5IHdpdGggYSBiaWcgcGFydHkgYW5kIGxvdHMgb2YgZnJpZW5kcy4KMTAuIFRoZSBwcm9qZWN0IGlzIGR1ZSBpbiB0ZW4gZGF5cywgc28gd2UgbmVlZCB0byB3b3JrIHF1aWNrbHku"
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("C10").Value = $text

$b64 = "U3VyZSEgQmVsb3cgaXMgYSBzaW1wbGUgaW1wbGVtZW50YXRpb24gb2YgdGhlIGNsYXNzaWMgU25ha2UgZ2FtZSB1c2luZyBQeXRob24gYW5kIHRoZSBQeWdhbWUgbGlicmFyeS4gVG8gcnVuIHRoaXMgY29kZSwgeW91IG5lZWQgdG8gaGF2ZSBQeXRob24gYW5kIFB5Z2FtZSBpbnN0YWxsZWQgb24geW91ciBzeXN0ZW0uCgpGaXJzdCwgeW91IGNhbiBpbnN0YWxsIFB5Z2FtZSB1c2luZyBwaXAgaWYgeW91IGhhdmVuJ3QgYWxyZWFkeToKYGBgc2gKcGlwIGluc3RhbGwgcHlnYW1lCmBgYAoKTm93LCBoZXJlIGlzIHRoZSBjb21wbGV0ZSBjb2RlIGZvciB0aGUgU25ha2UgZ2FtZToKCmBgYHB5dGhvbgppbXBvcnQgcHlnYW1lCmltcG9ydCB0aW1lCmltcG9ydCByYW5kb20KCiMgSW5pdGlhbGl6ZSBQeWdhbWUKcHlnYW1lLmluaXQoKQoKIyBEZWZpbmUgY29sb3JzCndoaXRlID0gKDI1NSwgMjU1LCAyNTUpCnllbGxvdyA9ICgyNTUsIDI1NSwgMTAyKQpibGFjayA9ICgwLCAwLCAwKQpyZWQgPSAoMjEzLCA1MCwgODApCmdyZWVuID0gKDAsIDI1NSwgMCkKYmx1ZSA9ICg1MCwgMTUzLCAyMTMpCgojIERlZmluZSBkaXNwbGF5IGRpbWVuc2lvbnMKZGlzX3dpZHRoID0gODAwCmRpc19oZWlnaHQgPSA2MDAKCiMgQ3JlYXRlIHRoZSBkaXNwbGF5CmRpcyA9IHB5Z2FtZS5kaXNwbGF5LnNldF9tb2RlKChkaXNfd2lkdGgsIGRpc19oZWlnaHQpKQpweWdhbWUuZGlzcGxheS5zZXRfY2FwdGlvbignU25ha2UgR2FtZSBieSBZb3VyTmFtZScpCgojIERlZmluZSB0aGUgY2xvY2sKY2xvY2sgPSBweWdhbWUudGltZS5DbG9jaygpCgojIERlZmluZSBzbmFrZSBibG9jayBzaXplIGFuZCBzcGVlZApzbmFrZV9ibG9jayA9IDEwCnNuYWtlX3NwZWVkID0gMTUKCiMgRGVmaW5lIGZvbnRzCmZvbnRfc3R5bGUgPSBweWdhbWUuZm9udC5TeXNGb250KCJiYWhuc2NocmlmdCIsIDI1KQpzY29yZV9mb250ID0gcHlnYW1lLmZvbnQuU3lzRm9udCgiY29taWNzYW5zbXMiLCAzNSkKCmRlZiBvdXJfc25ha2Uoc25ha2VfYmxvY2ssIHNuYWtlX2xpc3QpOgogICAgZm9yIHggaW4gc25ha2VfbGlzdDoKICAgICAgICBweWdhbWUuZHJhdy5yZWN0KGRpcywgYmxhY2ssIFt4WzBdLCB4WzFdLCBzbmFrZV9ibG9jaywgc25ha2VfYmxvY2tdKQoKZGVmIG1lc3NhZ2UobXNnLCBjb2xvcik6CiAgICBtZXNnID0gZm9udF9zdHlsZS5yZW5kZXIobXNnLCBUcnVlLCBjb2xvcikKICAgIGRpcy5ibGl0KG1lc2csIFtkaXNfd2lkdGggLyA2LCBkaXNfaGVpZ2h0IC8gM10pCgpkZWYgZ2FtZUxvb3AoKToKICAgIGdhbWVfb3ZlciA9IEZhbHNlCiAgICBnYW1lX2Nsb3NlID0gRmFsc2UKCiAgICB4MSA9IGRpc193aWR0aCAvIDIKICAgIHkxID0gZGlzX2hlaWdodCAvIDIKCiAgICB4MV9jaGFuZ2UgPSAwCiAgICB5MV9jaGFuZ2UgPSAwCgogICAgc25ha2VfTGlzdCA9IFtdCiAgICBMZW5ndGhfb2Zfc25ha2UgPSAxCgogICAgZm9vZHggPSByb3VuZChyYW5kb20ucmFuZHJhbmdlKDAsIGRpc193aWR0aCAtIHNuYWtlX2Jsb2NrKSAvIDEwLjApICogMTAuMAogICAgZm9vZHkgPSByb3VuZChyYW5kb20ucmFuZHJhbmdlKDAsIGRpc19oZWlnaHQgLSBzbmFrZV9ibG9jaykgLyAxMC4wKSAqIDEwLjAKCiAgICB3aGlsZSBub3QgZ2FtZV9vdmVyOgoKICAgICAgICB3aGlsZSBnYW1lX2Nsb3NlID09IFRydWU6CiAgICAgICAgICAgIGRpcy5maWxsKGJsdWUpCiAgICAgICAgICAgIG1lc3NhZ2UoIllvdSBMb3N0ISBQcmVzcyBRLVF1aXQgb3IgQy1QbGF5IEFnYWluIiwgcmVkKQogICAgICAgICAgICBweWdhbWUuZGlzcGxheS51cGRhdGUoKQoKICAgICAgICAgICAgZm9yIGV2ZW50IGluIHB5Z2FtZS5ldmVudC5nZXQoKToKICAgICAgICAgICAgICAgIGlmIGV2ZW50LnR5cGUgPT0gcHlnYW1lLktFWURPV046CiAgICAgICAgICAgICAgICAgICAgaWYgZXZlbnQua2V5ID09IHB5Z2FtZS5LX3E6CiAgICAgICAgICAgICAgICAgICAgICAgIGdhbWVfb3ZlciA9IFRydWUKICAgICAgICAgICAgICAgICAgICAgICAgZ2FtZV9jbG9zZSA9IEZhbHNlCiAgICAgICAgICAgICAgICAgICAgaWYgZXZlbnQua2V5ID09IHB5Z2FtZS5LX2M6CiAgICAgICAgICAgICAgICAgICAgICAgIGdhbWVMb29wKCkKCiAgICAgICAgZm9yIGV2ZW50IGluIHB5Z2FtZS5ldmVudC5nZXQoKToKICAgICAgICAgICAgaWYgZXZlbnQudHlwZSA9PSBweWdhbWUuUVVJVDoKICAgICAgICAgICAgICAgIGdhbWVfb3ZlciA9IFRydWUKICAgICAgICAgICAgaWYgZXZlbnQudHlwZSA9PSBweWdhbWUuS0VZRE9XTjoKICAgICAgICAgICAgICAgIGlmIGV2ZW50LmtleSA9PSBweWdhbWUuS19MRUZUOgogICAgICAgICAgICAgICAgICAgIHgxX2NoYW5nZSA9IC1zbmFrZV9ibG9jawogICAgICAgICAgICAgICAgICAgIHkxX2NoYW5nZSA9IDAKICAgICAgICAgICAgICAgIGVsaWYgZXZlbnQua2V5ID09IHB5Z2FtZS5LX1JJR0hUOgogICAgICAgICAgICAgICAgICAgIHgxX2NoYW5nZSA9IHNuYWtlX2Jsb2NrCiAgICAgICAgICAgICAgICAgICAgeTFfY2hhbmdlID0gMAogICAgICAgICAgICAgICAgZWxpZiBldmVudC5rZXkgPT0gcHlnYW1lLktfVVA6CiAgICAgICAgICAgICAgICAgICAgeTFfY2hhbmdlID0gLXNuYWtlX2Jsb2NrCiAgICAgICAgICAgICAgICAgICAgeDFfY2hhbmdlID0gMAogICAgICAgICAgICAgICAgZWxpZiBldmVudC5rZXkgPT0gcHlnYW1lLktfRE9XTjoKICAgICAgICAgICAgICAgICAgICB5MV9jaGFuZ2UgPSBzbmFrZV9ibG9jawogICAgICAgICAgICAgICAgICAgIHgxX2NoYW5nZSA9IDAKCiAgICAgICAgaWYgeDEgPj0gZGlzX3dpZHRoIG9yIHgxIDwgMCBvciB5MSA+PSBkaXNfaGVpZ2h0IG9yIHkxIDwgMDoKICAgICAgICAgICAgZ2FtZV9jbG9zZSA9IFRydWUKICAgICAgICB4MSArPSB4MV9jaGFuZ2UKICAgICAgICB5MSArPSB5MV9jaGFuZ2UKICAgICAgICBkaXMuZmlsbChibHVlKQogICAgICAgIHB5Z2FtZS5kcmF3LnJlY3QoZGlzLCBncmVlbiwgW2Zvb2R4LCBmb29keSwgc25ha2VfYmxvY2ssIHNuYWtlX2Jsb2NrXSkKICAgICAgICBzbmFrZV9IZWFkID0gW10KICAgICAgICBzbmFrZV9IZWFkLmFwcGVuZCh4MSkKICAgICAgICBzbmFrZV9IZWFkLmFwcGVuZCh5MSkKICAgICAgICBzbmFrZV9MaXN0LmFwcGVuZChzbmFrZV9IZWFkKQogICAgICAgIGlmIGxlbihzbmFrZV9MaXN0KSA+IExlbmd0aF9vZl9zbmFrZToKICAgICAgICAgICAgZGVsIHNuYWtlX0xpc3RbMF0KCiAgICAgICAgZm9yIHggaW4gc25ha2VfTGlzdFs6LTFdOgogICAgICAgICAgICBpZiB4ID09IHNuYWtlX0hlYWQ6CiAgICAgICAgICAgICAgICBnYW1lX2Nsb3NlID0gVHJ1ZQoKICAgICAgICBvdXJfc25ha2Uoc25ha2VfYmxvY2ssIHNuYWtlX0xpc3QpCgogICAgICAgIHB5Z2FtZS5kaXNwbGF5LnVwZGF0ZSgpCgogICAgICAgIGlmIHgxID09IGZvb2R4IGFuZCB5MSA9PSBmb29keToKICAgICAgICAgICAgZm9vZHggPSByb3VuZChyYW5kb20ucmFuZHJhbmdlKDAsIGRpc193aWR0aCAtIHNuYWtlX2Jsb2NrKSAvIDEwLjApICogMTAuMAogICAgICAgICAgICBmb29keSA9IHJvdW5kKHJhbmRvbS5yYW5kcmFuZ2UoMCwgZGlzX2hlaWdodCAtIHNuYWtlX2Jsb2NrKSAvIDEwLjApICogMTAuMAogICAgICAgICAgICBMZW5ndGhfb2Zfc25ha2UgKz0gMQoKICAgICAgICBjbG9jay50aWNrKHNuYWtlX3NwZWVkKQoKICAgIHB5Z2FtZS5xdWl0KCkKICAgIHF1aXQoKQoKZ2FtZUxvb3AoKQpgYGAKClRoaXMgY29kZSBzZXRzIHVwIGEgYmFzaWMgU25ha2UgZ2FtZSB3aGVyZSB0aGUgcGxheWVyIGNvbnRyb2xzIGEgc25ha2UgdGhhdCBncm93cyBsb25nZXIgZWFjaCB0aW1lIGl0IGVhdHMgZm9vZC4gVGhlIGdhbWUgZW5kcyBpZiB0aGUgc25ha2UgcnVucyBpbnRvIHRoZSBzY3JlZW4gZWRnZXMgb3IgaXRzZWxmLiBUaGUgcGxheWVyIGNhbiByZXN0YXJ0IHRoZSBnYW1lIGJ5IHByZXNzaW5nICdDJyBvciBxdWl0IGJ5IHByZXNzaW5nICdRJyB3aGVuIHRoZSBnYW1lIGlzIG92ZXIuCgpUbyBydW4gdGhlIGdhbWUsIHNpbXBseSBzYXZlIHRoZSBjb2RlIHRvIGEgZmlsZSwgZm9yIGV4YW1wbGUgYHNuYWtlX2dhbWUucHlgLCBhbmQgcnVuIGl0IHVzaW5nIFB5dGhvbjoKCmBgYHNoCnB5dGhvbiBzbmFrZV9nYW1lLnB5CmBgYAoKRW5qb3kgeW91ciBnYW1lIQ=="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("C11").Value = $text

$b64 = "VGhlIG1vZGVsJ3MgYW5zd2VyIGlzIGNvcnJlY3QgYnV0IGluY2x1ZGVzIGV4dHJhIGluZm9ybWF0aW9uIGFuZCB1c2VzIGRpZmZlcmVudCB2YXJpYWJsZSBuYW1lcyBhbmQgY29sb3JzIGNvbXBhcmVkIHRvIHRoZSBleHBlY3RlZCBvdXRwdXQu"
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("E11").Value = $text

$b64 = "VG8gc29sdmUgdGhpcyBwcm9ibGVtLCB3ZSBuZWVkIHRvIGVuc3VyZSB0aGF0IGF0IG5vIHBvaW50IGFyZSB0aGUgbGlvbiBhbmQgZ29hdCBsZWZ0IGFsb25lIHRvZ2V0aGVyLCBub3IgdGhlIGxpb24gYW5kIGNhYmJhZ2UgbGVmdCBhbG9uZSB0b2dldGhlci4gSGVyZSBpcyBhIHN0ZXAtYnktc3RlcCBzb2x1dGlvbjoKCjEuICoqVGFrZSB0aGUgbGlvbiBhY3Jvc3MgdGhlIHJpdmVyIGZpcnN0LioqCiAgIC0gTm93IHRoZSBsaW9uIGlzIG9uIHRoZSBvdGhlciBzaWRlLCBhbmQgdGhlIGdvYXQgYW5kIGNhYmJhZ2UgYXJlIG9uIHRoZSBzdGFydGluZyBzaWRlLgoKMi4gKipSZXR1cm4gYWxvbmUgdG8gdGhlIHN0YXJ0aW5nIHNpZGUuKioKICAgLSBUaGUgbGlvbiBpcyBhbG9uZSBvbiB0aGUgb3RoZXIgc2lkZSwgYW5kIHRoZSBnb2F0IGFuZCBjYWJiYWdlIGFyZSBvbiB0aGUgc3RhcnRpbmcgc2lkZS4KCjMuICoqVGFrZSB0aGUgZ29hdCBhY3Jvc3MgdGhlIHJpdmVyLioqCiAgIC0gTm93IHRoZSBsaW9uIGFuZCBnb2F0IGFyZSBvbiB0aGUgb3RoZXIgc2lkZSwgYW5kIHRoZSBjYWJiYWdlIGlzIG9uIHRoZSBzdGFydGluZyBzaWRlLgoKNC4gKipSZXR1cm4gd2l0aCB0aGUgbGlvbiB0byB0aGUgc3RhcnRpbmcgc2lkZS4qKgogICAtIE5vdyB0aGUgZ29hdCBpcyBhbG9uZSBvbiB0aGUgb3RoZXIgc2lkZSwgYW5kIHRoZSBsaW9uIGFuZCBjYWJiYWdlIGFyZSBvbiB0aGUgc3RhcnRpbmcgc2lkZS4KCjUuICoqVGFrZSB0aGUgY2FiYmFnZSBhY3Jvc3MgdGhlIHJpdmVyLioqCiAgIC0gTm93IHRoZSBnb2F0IGFuZCBjYWJiYWdlIGFyZSBvbiB0aGUgb3RoZXIgc2lkZSwgYW5kIHRoZSBsaW9uIGlzIG9uIHRoZSBzdGFydGluZyBzaWRlLgoKNi4gKipSZXR1cm4gYWxvbmUgdG8gdGhlIHN0YXJ0aW5nIHNpZGUuKioKICAgLSBUaGUgZ29hdCBhbmQgY2FiYmFnZSBhcmUgb24gdGhlIG90aGVyIHNpZGUsIGFuZCB0aGUgbGlvbiBpcyBvbiB0aGUgc3RhcnRpbmcgc2lkZS4KCjcuICoqVGFrZSB0aGUgbGlvbiBhY3Jvc3MgdGhlIHJpdmVyIGFnYWluLioqCiAgIC0gTm93IGFsbCB0aHJlZSAobGlvbiwgZ29hdCwgYW5kIGNhYmJhZ2UpIGFyZSBvbiB0aGUgb3RoZXIgc2lkZS4KCkJ5IGZvbGxvd2luZyB0aGVzZSBzdGVwcywgeW91IGVuc3VyZSB0aGF0IHRoZSBsaW9uIGlzIG5ldmVyIGxlZnQgYWxvbmUgd2l0aCB0aGUgZ29hdCwgYW5kIHRoZSBsaW9uIGlzIG5ldmVyIGxlZnQgYWxvbmUgd2l0aCB0aGUgY2FiYmFnZS4="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("C12").Value = $text

$b64 = "VGhlIG1vZGVsJ3MgYW5zd2VyIG1hdGNoZXMgdGhlIGV4cGVjdGVkIG91dHB1dCBpbiBsb2dpYyBhbmQgc3RlcHMu"
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("E12").Value = $text

$b64 = "VG8gY2FsY3VsYXRlIHRoZSBSZXR1cm4gb24gSW52ZXN0bWVudCAoUk9JKSBmb3IgVmVnYW4gU3RlYWtzIGZvciB0aGUgeWVhciwgd2UgbmVlZCB0byBmb2xsb3cgdGhlc2Ugc3RlcHM6CgoxLiAqKkRldGVybWluZSB0aGUgTmV0IEJvb2sgVmFsdWUgb2YgQXNzZXRzIGF0IHRoZSBCZWdpbm5pbmcgb2YgdGhlIFllYXI6KioKICAgLSBDb3N0IG9mIGFzc2V0czogJDIwLDAwMCwwMDAKICAgLSBBY2N1bXVsYXRlZCBkZXByZWNpYXRpb246ICQ1LDAwMCwwMDAKICAgLSBOZXQgYm9vayB2YWx1ZSBhdCB0aGUgYmVnaW5uaW5nIG9mIHRoZSB5ZWFyID0gQ29zdCBvZiBhc3NldHMgLSBBY2N1bXVsYXRlZCBkZXByZWNpYXRpb24KICAgLSBOZXQgYm9vayB2YWx1ZSBhdCB0aGUgYmVnaW5uaW5nIG9mIHRoZSB5ZWFyID0gJDIwLDAwMCwwMDAgLSAkNSwwMDAsMDAwID0gJDE1LDAwMCwwMDAKCjIuICoqQ2FsY3VsYXRlIHRoZSBOZXQgQm9vayBWYWx1ZSBvZiBBc3NldHMgYXQgdGhlIEVuZCBvZiB0aGUgWWVhcjoqKgogICAtIERlcHJlY2lhdGlvbiBleHBlbnNlIGZvciB0aGUgeWVhcjogJDEsMDAwLDAwMAogICAtIEFjY3VtdWxhdGVkIGRlcHJlY2lhdGlvbiBhdCB0aGUgZW5kIG9mIHRoZSB5ZWFyID0gQWNjdW11bGF0ZWQgZGVwcmVjaWF0aW9uIGF0IHRoZSBiZWdpbm5pbmcgb2YgdGhlIHllYXIgKyBEZXByZWNpYXRpb24gZXhwZW5zZSBmb3IgdGhlIHllYXIKICAgLSBBY2N1bXVsYXRlZCBkZXByZWNpYXRpb24gYXQgdGhlIGVuZCBvZiB0aGUgeWVhciA9ICQ1LDAwMCwwMDAgKyAkMSwwMDAsMDAwID0gJDYsMDAwLDAwMAogICAtIE5ldCBib29rIHZhbHVlIGF0IHRoZSBlbmQgb2YgdGhlIHllYXIgPSBDb3N0IG9mIGFzc2V0cyAtIEFjY3VtdWxhdGVkIGRlcHJlY2lhdGlvbiBhdCB0aGUgZW5kIG9mIHRoZSB5ZWFyCiAgIC0gTmV0IGJvb2sgdmFsdWUgYXQgdGhlIGVuZCBvZiB0aGUgeWVhciA9ICQyMCwwMDAsMDAwIC0gJDYsMDAwLDAwMCA9ICQxNCwwMDAsMDAwCgozLiAqKkNhbGN1bGF0ZSB0aGUgQXZlcmFnZSBOZXQgQm9vayBWYWx1ZSBvZiBBc3NldHMgZm9yIHRoZSBZZWFyOioqCiAgIC0gQXZlcmFnZSBuZXQgYm9vayB2YWx1ZSA9IChOZXQgYm9vayB2YWx1ZSBhdCB0aGUgYmVnaW5uaW5nIG9mIHRoZSB5ZWFyICsgTmV0IGJvb2sgdmFsdWUgYXQgdGhlIGVuZCBvZiB0aGUgeWVhcikgLyAyCiAgIC0gQXZlcmFnZSBuZXQgYm9vayB2YWx1ZSA9ICgkMTUsMDAwLDAwMCArICQxNCwwMDAsMDAwKSAvIDIgPSAkMTQsNTAwLDAwMAoKNC4gKipDYWxjdWxhdGUgdGhlIFJPSToqKgogICAtIE9wZXJhdGluZyBwcm9maXQgZm9yIHRoZSB5ZWFyOiAkOTUwLDAwMAogICAtIFJPSSA9IChPcGVyYXRpbmcgcHJvZml0IC8gQXZlcmFnZSBuZXQgYm9vayB2YWx1ZSBvZiBhc3NldHMpICogMTAwCiAgIC0gUk9JID0gKCQ5NTAsMDAwIC8gJDE0LDUwMCwwMDApICogMTAwCiAgIC0gUk9JID0gMC4wNjU1ICogMTAwCiAgIC0gUk9JID0gNi41NSUKClRoZXJlZm9yZSwgdGhlIFJPSSBmb3IgVmVnYW4gU3RlYWtzIGZvciB0aGUgeWVhciBpcyAqKjYuNTUlKiou"
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("C13").Value = $text

$b64 = "TGV0J3MgYnJlYWsgZG93biB0aGUgc2NlbmFyaW8gc3RlcCBieSBzdGVwIHRvIGVuc3VyZSB3ZSB1bmRlcnN0YW5kIGl0IGNvcnJlY3RseToKCjEuIEluaXRpYWxseSwgdGhlcmUgYXJlIDIga2lsbGVycyBpbiB0aGUgcm9vbS4KMi4gT25lIGd1eSBlbnRlcnMgdGhlIHJvb20uCjMuIFRoaXMgZ3V5IGtpbGxzIDEgb2YgdGhlIDIga2lsbGVycy4KCk5vdywgbGV0J3MgYW5hbHl6ZSB0aGUgc2l0dWF0aW9uOgoKLSBJbml0aWFsbHk6IDIga2lsbGVycyBpbiB0aGUgcm9vbS4KLSBBZnRlciB0aGUgZ3V5IGVudGVycyBhbmQga2lsbHMgMSBraWxsZXI6IDEga2lsbGVyIGlzIGRlYWQsIGxlYXZpbmcgMSBraWxsZXIgcmVtYWluaW5nLgoKSG93ZXZlciwgd2UgbXVzdCBhbHNvIGNvbnNpZGVyIHRoZSBzdGF0dXMgb2YgdGhlIGd1eSB3aG8gZW50ZXJlZCB0aGUgcm9vbS4gQnkga2lsbGluZyBvbmUgb2YgdGhlIGtpbGxlcnMsIGhlIGhpbXNlbGYgYmVjb21lcyBhIGtpbGxlci4KClNvLCBhZnRlciB0aGUgZ3V5IGtpbGxzIG9uZSBvZiB0aGUgb3JpZ2luYWwga2lsbGVycywgdGhlIHJvb20gbm93IGNvbnRhaW5zOgotIFRoZSByZW1haW5pbmcgb3JpZ2luYWwga2lsbGVyICgxIGtpbGxlcikuCi0gVGhlIGd1eSB3aG8ga2lsbGVkIHRoZSBvcmlnaW5hbCBraWxsZXIgKG5vdyBhbHNvIGEga2lsbGVyKS4KClRoZXJlZm9yZSwgdGhlIHRvdGFsIG51bWJlciBvZiBraWxsZXJzIGxlZnQgaW4gdGhlIHJvb20gaXM6CjEgKHJlbWFpbmluZyBvcmlnaW5hbCBraWxsZXIpICsgMSAodGhlIGd1eSB3aG8gYmVjYW1lIGEga2lsbGVyKSA9IDIga2lsbGVycy4KClNvLCB0aGVyZSBhcmUgMiBraWxsZXJzIGxlZnQgaW4gdGhlIHJvb20u"
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("C14").Value = $text

$b64 = "VGhlIG1vZGVsJ3MgYW5zd2VyIGNvcnJlY3RseSBpZGVudGlmaWVzIHRoYXQgdGhlcmUgYXJlIDIga2lsbGVycyBsZWZ0IGluIHRoZSByb29tLg=="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("E14").Value = $text

$b64 = "SGVyZSBpcyB0aGUgZ2l2ZW4gYm9va3N0b3JlIGludmVudG9yeSBpbmZvcm1hdGlvbiBjb252ZXJ0ZWQgaW50byBKU09OIGZvcm1hdDoKCmBgYGpzb24KewogICJCb29rc3RvcmUgSW52ZW50b3J5IjogWwogICAgewogICAgICAiVGl0bGUiOiAiVG8gS2lsbCBhIE1vY2tpbmdiaXJkIiwKICAgICAgIkF1dGhvciI6ICJIYXJwZXIgTGVlIiwKICAgICAgIlF1YW50aXR5IjogMzAKICAgIH0sCiAgICB7CiAgICAgICJUaXRsZSI6ICIxOTg0IiwKICAgICAgIkF1dGhvciI6ICJHZW9yZ2UgT3J3ZWxsIiwKICAgICAgIlF1YW50aXR5IjogMTUKICAgIH0sCiAgICB7CiAgICAgICJUaXRsZSI6ICJUaGUgR3JlYXQgR2F0c2J5IiwKICAgICAgIkF1dGhvciI6ICJGLiBTY290dCBGaXR6Z2VyYWxkIiwKICAgICAgIlF1YW50aXR5IjogMjAKICAgIH0KICBdCn0KYGBg"
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("C15").Value = $text

$b64 = "VGhlIG1vZGVsJ3MgYW5zd2VyIGlzIGNvcnJlY3QgYW5kIG1hdGNoZXMgdGhlIGV4cGVjdGVkIG91dHB1dC4="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("E15").Value = $text

$b64 = "VG8gZGV0ZXJtaW5lIHRoZSBtYXhpbXVtIGFubnVhbCB2YWNhdGlvbiBhY2NydWFsIGZvciBmdWxsLXRpbWUgZW1wbG95ZWVzIHdpdGggbW9yZSB0aGFuIDEwIHllYXJzIG9mIGNyZWRpdGVkIHNlcnZpY2UsIHdlIG5lZWQgdG8gbG9vayBhdCB0aGUgcHJvdmlkZWQgZGF0YSBmb3IgZnVsbC10aW1lIGVtcGxveWVlcyBhbmQgdGhlIGNvcnJlc3BvbmRpbmcgeWVhcnMgb2Ygc2VydmljZS4KClRoZSByZWxldmFudCBzZWN0aW9uIGZvciBmdWxsLXRpbWUgZW1wbG95ZWVzIGlzOgoKLSBGdWxsLXRpbWUgRW1wbG95ZWVzICg4MCBob3VycyBwZXIgYmktd2Vla2x5IHBheSBwZXJpb2QpCiAgLSBNYXhpbXVtIEFubnVhbCBWYWNhdGlvbiBBY2NydWFsCiAgICAtIDEyMCBIb3VycyAoMTUgRGF5cykKICAgIC0gMTYwIEhvdXJzICgyMCBEYXlzKQogICAgLSAyMDAgSG91cnMgKDI1IERheXMpCiAgLSBNYXhpbXVtIFZhY2F0aW9uIEFjY3J1YWwgQ2FwCiAgICAtIDE4MCBob3VycwogICAgLSAyNDAgaG91cnMKICAgIC0gMzAwIGhvdXJzCiAgLSBZZWFycyBvZiBDcmVkaXRlZCBTZXJ2aWNlcwogICAgLSA8IDUgWWVhcnMKICAgIC0gNS0xMCBZZWFycwogICAgLSA+IDEwIFllYXJzCgpGb3IgZW1wbG95ZWVzIHdpdGggbW9yZSB0aGFuIDEwIHllYXJzIG9mIGNyZWRpdGVkIHNlcnZpY2UsIHRoZSBtYXhpbXVtIGFubnVhbCB2YWNhdGlvbiBhY2NydWFsIGlzOgotIDIwMCBIb3VycyAoMjUgRGF5cykKClRoZXJlZm9yZSwgdGhlIG1heGltdW0gYW5udWFsIHZhY2F0aW9uIGFjY3J1YWwgZm9yIGZ1bGwtdGltZSBlbXBsb3llZXMgd2l0aCBtb3JlIHRoYW4gMTAgeWVhcnMgb2YgY3JlZGl0ZWQgc2VydmljZSBpcyAyMDAgaG91cnMgKDI1IGRheXMpLg=="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("C16").Value = $text

$b64 = "VGhlIG1vZGVsJ3MgYW5zd2VyIGNvcnJlY3RseSBpZGVudGlmaWVzIHRoZSBtYXhpbXVtIGFubnVhbCB2YWNhdGlvbiBhY2NydWFsIGZvciBmdWxsLXRpbWUgZW1wbG95ZWVzIHdpdGggbW9yZSB0aGFuIDEwIHllYXJzIG9mIGNyZWRpdGVkIHNlcnZpY2UgYXMgMjAwIGhvdXJzICgyNSBkYXlzKS4="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("E16").Value = $text

$b64 = "QmFzZWQgb24gdGhlIHByb3ZpZGVkIGluZm9ybWF0aW9uIGZyb20gdGhlIFJlZGhvcnNlIFBhaWQgVGltZSBPZmYgUHJvZ3JhbSBGQVFzIGVmZmVjdGl2ZSAxLjEuMjAyMiwgaGVyZSdzIHdoYXQgaGFwcGVucyBpZiB5b3UgdXNlZCB2YWNhdGlvbiB0aGF0IHdhcyBub3QgYWNjcnVlZCBhbmQgeW91IGxlYXZlIHRoZSBjb21wYW55OgoKMS4gKipOb24tQ2FsaWZvcm5pYSBFbXBsb3llZXM6KioKICAgLSBJZiB5b3UgaGF2ZSBhIG5lZ2F0aXZlIHZhY2F0aW9uIGJhbGFuY2UgKGkuZS4sIHlvdSB1c2VkIG1vcmUgdmFjYXRpb24gdGhhbiB5b3UgaGFkIGFjY3J1ZWQpIGFuZCB5b3VyIGVtcGxveW1lbnQgaXMgdGVybWluYXRlZCBmb3IgYW55IHJlYXNvbiwgdGhlIGNvbXBhbnkgaXMgYXV0aG9yaXplZCB0byB3aXRoaG9sZCB0aGUgYW1vdW50IG93ZWQgZnJvbSB5b3VyIGZpbmFsIHBheWNoZWNrLgogICAtIElmIHRoZSBhbW91bnQgb3dlZCBleGNlZWRzIHlvdXIgZmluYWwgcGF5Y2hlY2ssIHlvdSB3aWxsIG5lZWQgdG8gcGF5IHRoZSByZW1haW5pbmcgYmFsYW5jZSB0byBSZWRob3JzZSB3aXRoaW4gMzAgZGF5cyBvZiB5b3VyIGxhc3QgZGF5IG9mIHdvcmsuCgoyLiAqKkNhbGlmb3JuaWEgRW1wbG95ZWVzOioqCiAgIC0gSWYgeW91IGhhdmUgYSBuZWdhdGl2ZSB2YWNhdGlvbiBiYWxhbmNlIGFuZCB5b3VyIGVtcGxveW1lbnQgaXMgdGVybWluYXRlZCBmb3IgYW55IHJlYXNvbiwgeW91IGFyZSByZXF1aXJlZCB0byBwYXkgYmFjayB0aGUgYW1vdW50IG93ZWQgdG8gUmVkaG9yc2Ugd2l0aGluIDMwIGRheXMgb2YgeW91ciBsYXN0IGRheSBvZiB3b3JrLgogICAtIFRoZSBwb2xpY3kgZG9lcyBub3Qgc3BlY2lmeSB3aXRoaG9sZGluZyBmcm9tIHRoZSBmaW5hbCBwYXljaGVjayBmb3IgQ2FsaWZvcm5pYSBlbXBsb3llZXMsIGJ1dCBpdCBkb2VzIHN0YXRlIHRoYXQgYWxsIG1vbmllcyBvd2VkIG11c3QgYmUgcGFpZCBiYWNrIHdpdGhpbiB0aGUgc3BlY2lmaWVkIHRpbWVmcmFtZS4KCkluIHN1bW1hcnksIGlmIHlvdSBsZWF2ZSB0aGUgY29tcGFueSB3aXRoIGEgbmVnYXRpdmUgdmFjYXRpb24gYmFsYW5jZSwgeW91IHdpbGwgbmVlZCB0byByZXBheSB0aGUgYW1vdW50IG93ZWQgdG8gUmVkaG9yc2UsIGVpdGhlciB0aHJvdWdoIHdpdGhob2xkaW5nIGZyb20geW91ciBmaW5hbCBwYXljaGVjayAoZm9yIG5vbi1DYWxpZm9ybmlhIGVtcGxveWVlcykgb3IgYnkgZGlyZWN0IHBheW1lbnQgd2l0aGluIDMwIGRheXMgKGZvciBib3RoIG5vbi1DYWxpZm9ybmlhIGFuZCBDYWxpZm9ybmlhIGVtcGxveWVlcyku"
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("C17").Value = $text

$b64 = "VGhlIG1vZGVsJ3MgYW5zd2VyIGlzIGNvcnJlY3QgYnV0IGluY2x1ZGVzIGV4dHJhIGRldGFpbHMgbm90IHByZXNlbnQgaW4gdGhlIGV4cGVjdGVkIG91dHB1dC4="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("E17").Value = $text

